$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Column A (test case), Column B (result), Column C (browser)
$newRows = @(
    ,@(91, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(92, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(93, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(94, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(95, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(96, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(97, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(98, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(99, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(100, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(101, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(102, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(103, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(104, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(105, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "PASSED", "chrome")
    ,@(106, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(107, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(108, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "PASSED", "chrome")
    ,@(109, "Edit", "FAILED", "chrome")
    ,@(110, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "PASSED", "chrome")
    ,@(111, "Edit", "PASSED", "chrome")
    ,@(112, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(113, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(114, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(115, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(116, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "PASSED", "chrome")
    ,@(117, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(118, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "PASSED", "chrome")
    ,@(119, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "FAILED", "chrome")
    ,@(120, "The User Add Edit And Delete Document Types under Setup -> Parameters Document Types", "PASSED", "chrome")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
}

